$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 341 (shifts existing rows 341..404 down to 342..405)
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row 341 with the new data record
$ws.Cells.Item(341, 1).Value = 5
$ws.Cells.Item(341, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(341, 3).Value = "Maule"
$ws.Cells.Item(341, 4).Value = 45015
$ws.Cells.Item(341, 5).Value = 7
$ws.Cells.Item(341, 6).Value = 100112009
$ws.Cells.Item(341, 7).Value = "Acelga"
$ws.Cells.Item(341, 8).Value = "Sin especificar"
$ws.Cells.Item(341, 9).Value = "Primera"
$ws.Cells.Item(341, 10).Value = 500
$ws.Cells.Item(341, 11).Value = 2500
$ws.Cells.Item(341, 12).Value = 2500
$ws.Cells.Item(341, 13).Value = 2500
$ws.Cells.Item(341, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(341, 15).Value = "Región del Maule"
$ws.Cells.Item(341, 16).Value = 625
$ws.Cells.Item(341, 17).Value = 4
$ws.Cells.Item(341, 18).Value = "Hortaliza"
